$d = $word.ActiveDocument

# 1. Refresh the timestamp in the footer ("2025-06-30 12:13Z / " -> "2025-07-02 02:48Z / ")
$f = $d.Sections(1).Footers(1)
$f.Range.Find.Execute("2025-06-30 12:13Z / ", $true, $false, $false, $false, $false, `
                       $true, 1, $false, "2025-07-02 02:48Z / ", 2)

# 2. Add the regression-test character styles: b, i, sub, sup, u
#    Each is based on "Default Paragraph Font", uiPriority 1, quick style.

$sb = $d.Styles.Add("b", 2)
$sb.BaseStyle = "DefaultParagraphFont"
$sb.Priority = 1
$sb.QuickStyle = $true
$sb.Font.Bold = $true

$si = $d.Styles.Add("i", 2)
$si.BaseStyle = "DefaultParagraphFont"
$si.Priority = 1
$si.QuickStyle = $true
$si.Font.Italic = $true

$ssub = $d.Styles.Add("sub", 2)
$ssub.BaseStyle = "DefaultParagraphFont"
$ssub.Priority = 1
$ssub.QuickStyle = $true
$ssub.Font.Subscript = $true

$ssup = $d.Styles.Add("sup", 2)
$ssup.BaseStyle = "DefaultParagraphFont"
$ssup.Priority = 1
$ssup.QuickStyle = $true
$ssup.Font.Superscript = $true

$su = $d.Styles.Add("u", 2)
$su.BaseStyle = "DefaultParagraphFont"
$su.Priority = 1
$su.QuickStyle = $true
$su.Font.Underline = 1
